$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("leadlag")

# Row 4 - new trade row appended to the leadlag sheet
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "'2026-02-16"
$ws.Range("C4").Value = "21:20:31"
$ws.Range("D4").Value = "leadlag"
$ws.Range("E4").Value = "DOWN"
$ws.Range("F4").Value = 69476.05
$ws.Range("H4").Value = "OPEN"
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0.75
$ws.Range("L4").Value = "Binance leading with -0.099% move"
$ws.Range("N4").Value = 0
